# This script updates the generated statistics (column F) on three of the
# four worksheets, matching a refreshed "gh-pages" data output.
#
#   展览   (Exhibitions)   - sheet 1
#   演出   (Performances)  - sheet 2
#   本地生活(Local life)    - sheet 3 (unchanged)
#   全部类型(All types)     - sheet 4
#
$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1) ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 97
$ws1.Range("F4").Value  = 1443
$ws1.Range("F5").Value  = 184
$ws1.Range("F6").Value  = 37
$ws1.Range("F7").Value  = 36
$ws1.Range("F8").Value  = 9699
$ws1.Range("F10").Value = 109
$ws1.Range("F11").Value = 241
$ws1.Range("F12").Value = 186
$ws1.Range("F13").Value = 367
$ws1.Range("F14").Value = 6677
$ws1.Range("F15").Value = 1082
$ws1.Range("F18").Value = 183

# --- 演出 (sheet2) ----------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 43

# --- 全部类型 (sheet4) -------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 97
$ws4.Range("F4").Value  = 1443
$ws4.Range("F5").Value  = 184
$ws4.Range("F6").Value  = 37
$ws4.Range("F7").Value  = 36
$ws4.Range("F8").Value  = 43
$ws4.Range("F10").Value = 9699
$ws4.Range("F12").Value = 109
$ws4.Range("F13").Value = 241
$ws4.Range("F14").Value = 186
$ws4.Range("F15").Value = 367
$ws4.Range("F16").Value = 6677
$ws4.Range("F17").Value = 1082
$ws4.Range("F20").Value = 183
